# Update "想去人数" (interest count) figures on the "展览" (rId1), "演出" (rId2)
# and "全部类型" (rId4) worksheets to match the refreshed scrape output.
# "本地生活" (rId3) is untouched by this update.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$ws1.Range("F2").Value  = 2810
$ws1.Range("F4").Value  = 365
$ws1.Range("F5").Value  = 1567
$ws1.Range("F6").Value  = 1156
$ws1.Range("F8").Value  = 548
$ws1.Range("F9").Value  = 127
$ws1.Range("F11").Value = 9483
$ws1.Range("F12").Value = 407
$ws1.Range("F13").Value = 2508
$ws1.Range("F16").Value = 184
$ws1.Range("F18").Value = 668
$ws1.Range("F21").Value = 1002
$ws1.Range("F22").Value = 2958
$ws1.Range("F23").Value = 2234
$ws1.Range("F30").Value = 14
$ws1.Range("F31").Value = 172
$ws1.Range("F32").Value = 219
$ws1.Range("F37").Value = 506
$ws1.Range("F38").Value = 22
$ws1.Range("F39").Value = 111
$ws1.Range("F40").Value = 1371
$ws1.Range("F41").Value = 118
$ws1.Range("F42").Value = 1468
$ws1.Range("F43").Value = 23
$ws1.Range("F44").Value = 335
$ws1.Range("F45").Value = 23
$ws1.Range("F46").Value = 359
$ws1.Range("F47").Value = 727

# --- 演出 (sheet2) ---
$ws2.Range("F5").Value  = 27

# --- 全部类型 (sheet4) ---
$ws4.Range("F2").Value  = 2810
$ws4.Range("F3").Value  = 365
$ws4.Range("F4").Value  = 1567
$ws4.Range("F6").Value  = 1156
$ws4.Range("F7").Value  = 548
$ws4.Range("F8").Value  = 127
$ws4.Range("F9").Value  = 9483
$ws4.Range("F10").Value = 407
$ws4.Range("F15").Value = 184
$ws4.Range("F16").Value = 668
$ws4.Range("F18").Value = 1002
$ws4.Range("F19").Value = 2959
$ws4.Range("F20").Value = 2234
$ws4.Range("F25").Value = 14
$ws4.Range("F26").Value = 172
$ws4.Range("F27").Value = 219
$ws4.Range("F32").Value = 506
$ws4.Range("F33").Value = 27
$ws4.Range("F36").Value = 22
$ws4.Range("F37").Value = 111
$ws4.Range("F38").Value = 1372
$ws4.Range("F40").Value = 118
$ws4.Range("F41").Value = 1468
$ws4.Range("F42").Value = 23
$ws4.Range("F44").Value = 335
$ws4.Range("F45").Value = 23
$ws4.Range("F46").Value = 359
$ws4.Range("F47").Value = 727

$wb.Save()
